# LOM3219.xlsx content update
# - Remove the last (now unused) template row (old row 26)
# - Update text content of several cells (the "Ficha de disciplina" rows
#   got re-mapped during the page's content rebuild)
# - Adjust a handful of row heights so the autofit-style formatting still
#   lines up with the (now shorter) text in each row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had one trailing template row that is no longer needed.
$ws.Rows(26).Delete() | Out-Null

# --- Cell content changes -------------------------------------------------

$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

$ws.Range("A17").Value = "Avaliação:"

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("A19").Value = "Critério:"

$ws.Range("A20").Value = "Norma de recuperação:"

$ws.Range("A21").Value = "Bibliografia:"

$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"
$ws.Range("C23").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"

$ws.Range("B24").Value = "LOM3229 -  Métodos Experimentais da Física II  (Requisito)`n"
$ws.Range("C24").Value = "LOM3229 -  Métodos Experimentais da Física II  (Requisito)`n"

$ws.Range("B25").Value = "LOM3246 -  Técnicas de Caracterização de Materiais  (Requisito)`n"
$ws.Range("C25").Value = "LOM3246 -  Técnicas de Caracterização de Materiais  (Requisito)`n"

# --- Row height adjustments ------------------------------------------------

$ws.Rows(13).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(17).RowHeight = 15
$ws.Rows(17).UseStandardHeight = $true
$ws.Rows(18).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(22).RowHeight = 15
$ws.Rows(22).UseStandardHeight = $true
$ws.Rows(23).RowHeight = 30
